$d = $word.ActiveDocument

# The title paragraph reads:
#   "TS Jatai – TS 3.4 Sanskrit Corrections – Observed till 31st Jan 2024"
# We need to append a period after "2024" (end of that paragraph),
# matching formatting (Bold, Underline, Size 32) of the surrounding text.

$r = $d.Content
$found = $r.Find.Execute("2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$r.Collapse(0)
$insStart = $r.Start
$r.InsertAfter(".")
$newRange = $d.Range($insStart, $insStart + 1)
$newRange.Font.Bold = $true
$newRange.Font.BoldBi = $true
$newRange.Font.Size = 16
$newRange.Font.SizeBi = 16
$newRange.Font.Underline = 1
